$wb = $excel.ActiveWorkbook

# --- Sheet "DATA_RAW" (sheet1): add row 23 "FAVÖK" ---
$ws1 = $wb.Worksheets.Item("DATA_RAW")
$ws1.Range("A23").Value = "FAVÖK"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 0
$ws1.Range("D23").Value = 0
$ws1.Range("E23").Value = 5939000000
$ws1.Range("F23").Value = 0
$ws1.Range("G23").Value = 0
$ws1.Range("H23").Value = 5999000000
$ws1.Range("I23").Value = 8513000000
$ws1.Range("J23").Value = 4438000000
$ws1.Range("K23").Value = 0
$ws1.Range("L23").Value = 4482000000
$ws1.Range("M23").Value = 8165000000
$ws1.Range("N23").Value = 5173000000
$ws1.Range("O23").Value = 0
$ws1.Range("P23").Value = 7116000000
$ws1.Range("Q23").Value = 7994000000
$ws1.Range("R23").Value = 5709000000
$ws1.Range("S23").Value = 0
$ws1.Range("T23").Value = 7203000000
$ws1.Range("U23").Value = 7771000000
$ws1.Range("V23").Value = 6374000000
$ws1.Range("W23").Value = 0
$ws1.Range("X23").Value = 5308000000
$ws1.Range("Y23").Value = 7969000000
$ws1.Range("Z23").Value = 7612000000
$ws1.Range("AA23").Value = 0
$ws1.Range("AB23").Value = 6334000000
$ws1.Range("AC23").Value = 7776000000
$ws1.Range("AD23").Value = 6974000000
$ws1.Range("AE23").Value = 0
$ws1.Range("AF23").Value = 5844000000
$ws1.Range("AG23").Value = 6026000000
$ws1.Range("AH23").Value = 6594000000
$ws1.Range("AI23").Value = -2053000000
$ws1.Range("AJ23").Value = 5793000000
$ws1.Range("AK23").Value = 7905000000
$ws1.Range("AL23").Value = 5283000000
$ws1.Range("AM23").Value = 3080000000
$ws1.Range("AN23").Value = 6715000000
$ws1.Range("AO23").Value = 8679000000
$ws1.Range("AP23").Value = 6723000000
$ws1.Range("AQ23").Value = 7682000000
$ws1.Range("AR23").Value = 7708000000
$ws1.Range("AS23").Value = 10258000000
$ws1.Range("AT23").Value = 8292000000
$ws1.Range("AU23").Value = 10379000000
$ws1.Range("AV23").Value = 9955000000
$ws1.Range("AW23").Value = 13891000000
$ws1.Range("AX23").Value = 10341000000
$ws1.Range("AY23").Value = 12405000000
$ws1.Range("AZ23").Value = 12686000000
$ws1.Range("BA23").Value = 17897000000
$ws1.Range("BB23").Value = 12975000000
$ws1.Range("BC23").Value = 13407000000
$ws1.Range("BD23").Value = 15876000000
$ws1.Range("BE23").Value = 22247000000
$ws1.Range("BF23").Value = 17048000000
$ws1.Range("BG23").Value = 0
$ws1.Range("BH23").Value = 20238000000
$ws1.Range("BI23").Value = 20399000000
$ws1.Range("BJ23").Value = 20364000000
$ws1.Range("BK23").Value = 0
$ws1.Range("BL23").Value = 21518000000
$ws1.Range("BM23").Value = 27032000000
$ws1.Range("BN23").Value = 22352000000
$ws1.Range("BO23").Value = 0
$ws1.Range("BP23").Value = 26895000000
$ws1.Range("BQ23").Value = 31653000000
$ws1.Range("BR23").Value = 27581000000
$ws1.Range("BS23").Value = 0
$ws1.Range("BT23").Value = 30552000000
$ws1.Range("BU23").Value = 32000000000
$ws1.Range("BV23").Value = 0
$ws1.Range("BW23").Value = 37961000000

# --- Sheet "gelir tablosu (çeyreklik)" (sheet3): add row 9 "FAVÖK" ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A9").Value = "FAVÖK"
$ws3.Range("B9").Value = 0
$ws3.Range("C9").Value = 0
$ws3.Range("D9").Value = 0
$ws3.Range("E9").Value = 5939000000
$ws3.Range("F9").Value = 0
$ws3.Range("G9").Value = 0
$ws3.Range("H9").Value = 5999000000
$ws3.Range("I9").Value = 8513000000
$ws3.Range("J9").Value = 4438000000
$ws3.Range("K9").Value = 0
$ws3.Range("L9").Value = 4482000000
$ws3.Range("M9").Value = 8165000000
$ws3.Range("N9").Value = 5173000000
$ws3.Range("O9").Value = 0
$ws3.Range("P9").Value = 7116000000
$ws3.Range("Q9").Value = 7994000000
$ws3.Range("R9").Value = 5709000000
$ws3.Range("S9").Value = 0
$ws3.Range("T9").Value = 7203000000
$ws3.Range("U9").Value = 7771000000
$ws3.Range("V9").Value = 6374000000
$ws3.Range("W9").Value = 0
$ws3.Range("X9").Value = 5308000000
$ws3.Range("Y9").Value = 7969000000
$ws3.Range("Z9").Value = 7612000000
$ws3.Range("AA9").Value = 0
$ws3.Range("AB9").Value = 6334000000
$ws3.Range("AC9").Value = 7776000000
$ws3.Range("AD9").Value = 6974000000
$ws3.Range("AE9").Value = 0
$ws3.Range("AF9").Value = 5844000000
$ws3.Range("AG9").Value = 6026000000
$ws3.Range("AH9").Value = 6594000000
$ws3.Range("AI9").Value = -2053000000
$ws3.Range("AJ9").Value = 5793000000
$ws3.Range("AK9").Value = 7905000000
$ws3.Range("AL9").Value = 5283000000
$ws3.Range("AM9").Value = 3080000000
$ws3.Range("AN9").Value = 6715000000
$ws3.Range("AO9").Value = 8679000000
$ws3.Range("AP9").Value = 6723000000
$ws3.Range("AQ9").Value = 7682000000
$ws3.Range("AR9").Value = 7708000000
$ws3.Range("AS9").Value = 10258000000
$ws3.Range("AT9").Value = 8292000000
$ws3.Range("AU9").Value = 10379000000
$ws3.Range("AV9").Value = 9955000000
$ws3.Range("AW9").Value = 13891000000
$ws3.Range("AX9").Value = 10341000000
$ws3.Range("AY9").Value = 12405000000
$ws3.Range("AZ9").Value = 12686000000
$ws3.Range("BA9").Value = 17897000000
$ws3.Range("BB9").Value = 12975000000
$ws3.Range("BC9").Value = 13407000000
$ws3.Range("BD9").Value = 15876000000
$ws3.Range("BE9").Value = 22247000000
$ws3.Range("BF9").Value = 17048000000
$ws3.Range("BG9").Value = 0
$ws3.Range("BH9").Value = 20238000000
$ws3.Range("BI9").Value = 20399000000
$ws3.Range("BJ9").Value = 20364000000
$ws3.Range("BK9").Value = 0
$ws3.Range("BL9").Value = 21518000000
$ws3.Range("BM9").Value = 27032000000
$ws3.Range("BN9").Value = 22352000000
$ws3.Range("BO9").Value = 0
$ws3.Range("BP9").Value = 26895000000
$ws3.Range("BQ9").Value = 31653000000
$ws3.Range("BR9").Value = 27581000000
$ws3.Range("BS9").Value = 0
$ws3.Range("BT9").Value = 30552000000
$ws3.Range("BU9").Value = 32000000000
$ws3.Range("BV9").Value = 0
$ws3.Range("BW9").Value = 37961000000
